$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right score and Wrong score corrected
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): total right/wrong marks and summary text corrected
$ws.Range("B12").Value = 96
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "92 / 112"
